$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "30.354.30"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.010.00"
$ws.Range("E3").Value = "  +4.93%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "324.92"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.5133"
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("D8").Value = "0.4265"
$ws.Range("E8").Value = "  +5.66%  "
$ws.Range("D9").Value = "0.08704"
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("D10").Value = "1.137"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("D11").Value = "43.27"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("D13").Value = "2.015.13"
$ws.Range("E13").Value = "  +5.22%  "
$ws.Range("D14").Value = "6.568"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "7.468"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "94.37"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "0.00001114"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").Value = "0.06534"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.90"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "6.203"
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("D23").Value = "30.396.59"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").Value = "2.245.08"
$ws.Range("E26").Value = "  +5.05%  "
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "162.26"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "2.423"
$ws.Range("E29").Value = "  +5.41%  "
$ws.Range("D30").Value = "130.98"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "1.137"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "3.825"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "1.374"
$ws.Range("E35").Value = "  +15.05%  "
$ws.Range("D36").Value = "0.02525"
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("D37").Value = "0.06677"
$ws.Range("E37").Value = "  +4.19%  "
$ws.Range("D38").Value = "5.463"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "12.38"
$ws.Range("E39").Value = "  +8.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.140"
$ws.Range("E40").Value = "  +5.08%  "
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("D42").Value = "0.6646"
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("D43").Value = "1.238"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "13.68"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("D46").Value = "0.6168"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").Value = "2.187"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "3.664"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.260"
$ws.Range("E49").Value = "  +4.59%  "
$ws.Range("D50").Value = "124.37"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").Value = "80.59"
$ws.Range("E51").Value = "  +2.14%  "
